$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("Q2").Value = 5.511848134097666
$ws.Range("R2").Value = 49.606633206879
$ws.Range("S2").Value = 0.007694456167922165
$ws.Range("T2").Value = 0.007694456167922165

# Row 3
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("Q3").Value = 316.089371518151
$ws.Range("R3").Value = 2844.804343663359
$ws.Range("S3").Value = 0.4412559553748731
$ws.Range("T3").Value = 0.4412559553748731

# Row 4
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("Q4").Value = 4.193650615478222
$ws.Range("S4").Value = 0.005854272479816698
$ws.Range("T4").Value = 0.005854272479816698

# Row 5
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("Q5").Value = 240.4943596346426
$ws.Range("S5").Value = 0.3357264684768396
$ws.Range("T5").Value = 0.3357264684768396

# Row 6
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("S6").Value = 0.003590037617071002
$ws.Range("T6").Value = 0.003590037617071002

# Row 7
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("S7").Value = 0.2058788098834775
$ws.Range("T7").Value = 0.2058788098834775
